$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dates in B2 and B3 (format already set to yyyy-mm-dd)
$ws.Range("B2").Value = 45659
$ws.Range("B3").Value = 45659

# Update the active selection to C3
$ws.Range("C3").Select()
